# Auto-generated edits applying market-price refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 20822
$ws.Range("J7").Value = 51253
$ws.Range("L7").Value = 51253
$ws.Range("N7").Value = -51477

$ws.Range("H13").Value = 2974.25
$ws.Range("I13").Value = 847.5
$ws.Range("J13").Value = 5101
$ws.Range("K13").Value = 847.5
$ws.Range("L13").Value = 5101
$ws.Range("M13").Value = -678.5
$ws.Range("N13").Value = -5439

$ws.Range("H14").Value = 20822
$ws.Range("J14").Value = 51253
$ws.Range("L14").Value = 51253
$ws.Range("N14").Value = -51635

$ws.Range("H16").Value = 19999
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 19999
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 19999
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -20459

$ws.Range("H19").Value = 1355.2
$ws.Range("I19").Value = 1006.8333
$ws.Range("J19").Value = 1504.5
$ws.Range("K19").Value = 1006.8333
$ws.Range("L19").Value = 1504.5
$ws.Range("M19").Value = -831.8333
$ws.Range("N19").Value = -1854.5

$ws.Range("H39").Value = 125.72222
$ws.Range("I39").Value = 15.8125
$ws.Range("J39").Value = 1005
$ws.Range("K39").Value = 47.4375
$ws.Range("L39").Value = 3015
$ws.Range("M39").Value = 248.5625
$ws.Range("N39").Value = -3607

$ws.Range("H70").Value = 145314.28
$ws.Range("J70").Value = 145314.28
$ws.Range("L70").Value = 435942.84
$ws.Range("N70").Value = -436482.84

$ws.Range("H73").Value = 145314.28
$ws.Range("J73").Value = 145314.28
$ws.Range("L73").Value = 435942.84
$ws.Range("N73").Value = -437814.84

$ws.Range("H100").Value = 1710.9
$ws.Range("J100").Value = 2002
$ws.Range("L100").Value = 2002
$ws.Range("N100").Value = -3084

$ws.Range("H132").Value = 1672.3959
$ws.Range("J132").Value = 3254.8
$ws.Range("L132").Value = 9764.400000000001
$ws.Range("N132").Value = -14824.4

$ws.Range("H141").Value = 779.96875
$ws.Range("I141").Value = 779.96875
$ws.Range("K141").Value = 2339.90625
$ws.Range("M141").Value = 2840.09375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4100.8887
$ws.Range("I32").Value = 2361.6562
$ws.Range("K32").Value = 2361.6562
$ws.Range("M32").Value = -2074.6562

$ws.Range("H61").Value = 2841.5217
$ws.Range("I61").Value = 2542.9744
$ws.Range("J61").Value = 4504.857
$ws.Range("K61").Value = 2542.9744
$ws.Range("L61").Value = 4504.857
$ws.Range("M61").Value = -2330.9744
$ws.Range("N61").Value = -4928.857

$ws.Range("H74").Value = 7090.875
$ws.Range("I74").Value = 2174.5386
$ws.Range("J74").Value = 12901.091
$ws.Range("K74").Value = 2174.5386
$ws.Range("L74").Value = 12901.091
$ws.Range("M74").Value = -1300.5386
$ws.Range("N74").Value = -14649.091

$ws.Range("H77").Value = 7090.875
$ws.Range("I77").Value = 2174.5386
$ws.Range("J77").Value = 12901.091
$ws.Range("K77").Value = 10872.693
$ws.Range("L77").Value = 64505.455
$ws.Range("M77").Value = -6504.692999999999
$ws.Range("N77").Value = -73241.455

$ws.Range("H132").Value = 3937.8235
$ws.Range("I132").Value = 3478.2173
$ws.Range("K132").Value = 10434.6519
$ws.Range("M132").Value = -7904.651899999999

$ws.Range("H136").Value = 2841.5217
$ws.Range("I136").Value = 2542.9744
$ws.Range("J136").Value = 4504.857
$ws.Range("K136").Value = 7628.9232
$ws.Range("L136").Value = 13514.571
$ws.Range("M136").Value = -5078.9232
$ws.Range("N136").Value = -18614.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3584.8462
$ws.Range("I20").Value = 2272.8125
$ws.Range("K20").Value = 2272.8125
$ws.Range("M20").Value = -2025.8125

$ws.Range("H107").Value = 5940.7617
$ws.Range("I107").Value = 1285.3334
$ws.Range("K107").Value = 1285.3334
$ws.Range("M107").Value = 634.6666

$ws.Range("H134").Value = 1435.2368
$ws.Range("I134").Value = 1441.174
$ws.Range("J134").Value = 1376.7142
$ws.Range("K134").Value = 4323.522
$ws.Range("L134").Value = 4130.142599999999
$ws.Range("M134").Value = -1788.522
$ws.Range("N134").Value = -9200.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 13233.667
$ws.Range("J69").Value = 20201
$ws.Range("L69").Value = 20201
$ws.Range("N69").Value = -21699

$ws.Range("H72").Value = 13233.667
$ws.Range("J72").Value = 20201
$ws.Range("L72").Value = 60603
$ws.Range("N72").Value = -68091

$ws.Range("H132").Value = 169665.64
$ws.Range("I132").Value = 335736.34
$ws.Range("J132").Value = 3594.9333
$ws.Range("K132").Value = 1007209.02
$ws.Range("L132").Value = 10784.7999
$ws.Range("M132").Value = -1004679.02
$ws.Range("N132").Value = -15844.7999

$ws.Range("H134").Value = 20138.371
$ws.Range("I134").Value = 14807.218
$ws.Range("J134").Value = 50792.5
$ws.Range("K134").Value = 44421.654
$ws.Range("L134").Value = 152377.5
$ws.Range("M134").Value = -41886.654
$ws.Range("N134").Value = -157447.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2447.9092
$ws.Range("I12").Value = 48.666668
$ws.Range("J12").Value = 2826.7368
$ws.Range("K12").Value = 146.000004
$ws.Range("L12").Value = 8480.2104
$ws.Range("M12").Value = 26.99999600000001
$ws.Range("N12").Value = -8826.2104

$ws.Range("H128").Value = 267994
$ws.Range("I128").Value = 267994
$ws.Range("K128").Value = 803982
$ws.Range("M128").Value = -799002

$ws.Range("H132").Value = 1121.0488
$ws.Range("I132").Value = 1090.7354
$ws.Range("J132").Value = 1268.2858
$ws.Range("K132").Value = 9816.6186
$ws.Range("L132").Value = 11414.5722
$ws.Range("M132").Value = -7286.6186
$ws.Range("N132").Value = -16474.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 10500
$ws.Range("J27").Value = 10500
$ws.Range("L27").Value = 10500
$ws.Range("N27").Value = -10832

$ws.Range("H70").Value = 11163.654
$ws.Range("I70").Value = 11649.875
$ws.Range("J70").Value = 10947.556
$ws.Range("K70").Value = 11649.875
$ws.Range("L70").Value = 10947.556
$ws.Range("M70").Value = -11379.875
$ws.Range("N70").Value = -11487.556

$ws.Range("H73").Value = 11163.654
$ws.Range("I73").Value = 11649.875
$ws.Range("J73").Value = 10947.556
$ws.Range("K73").Value = 11649.875
$ws.Range("L73").Value = 10947.556
$ws.Range("M73").Value = -10713.875
$ws.Range("N73").Value = -12819.556

$ws.Range("H119").Value = 29499.166
$ws.Range("J119").Value = 29499.166
$ws.Range("L119").Value = 29499.166
$ws.Range("N119").Value = -39175.166

$ws.Range("H132").Value = 3337.2
$ws.Range("I132").Value = 3260.3
$ws.Range("K132").Value = 9780.900000000001
$ws.Range("M132").Value = -7250.900000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3569.077
$ws.Range("I16").Value = 2831.6365
$ws.Range("J16").Value = 7625
$ws.Range("K16").Value = 2831.6365
$ws.Range("L16").Value = 7625
$ws.Range("M16").Value = -2661.6365
$ws.Range("N16").Value = -7965

$ws.Range("H55").Value = 250.61539
$ws.Range("I55").Value = 176.9
$ws.Range("K55").Value = 176.9
$ws.Range("M55").Value = -3.900000000000006

$ws.Range("H100").Value = 9071.909
$ws.Range("I100").Value = 11724
$ws.Range("K100").Value = 11724
$ws.Range("M100").Value = -11183

$ws.Range("H132").Value = 2480.3137
$ws.Range("I132").Value = 2327
$ws.Range("K132").Value = 6981
$ws.Range("M132").Value = -4451

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 52218.09
$ws.Range("J41").Value = 60407.832
$ws.Range("L41").Value = 60407.832
$ws.Range("N41").Value = -61187.832

$ws.Range("H51").Value = 21000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H132").Value = 1709.6666
$ws.Range("J132").Value = 1289
$ws.Range("L132").Value = 3867
$ws.Range("N132").Value = -8927

$ws.Range("H136").Value = 1222.1846
$ws.Range("I136").Value = 1102.1111
$ws.Range("K136").Value = 3306.3333
$ws.Range("M136").Value = -756.3333000000002
